$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 64
$ws.Range("C3").Value = 55
$ws.Range("C4").Value = 62
$ws.Range("C5").Value = 58
$ws.Range("C6").Value = 82
$ws.Range("C7").Value = 71
$ws.Range("C9").Value = 78
$ws.Range("C10").Value = 66
$ws.Range("C11").Value = 75
$ws.Range("C13").Value = 74
$ws.Range("C14").Value = 72

$ws.Range("C15").Select()
